$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 254, shifting existing rows 254-367 down to 255-368.
$ws.Rows.Item(254).Insert()

# Populate the newly inserted row 254 with the new record.
$ws.Range("A254").Value = 10
$ws.Range("B254").Value = "Vega Modelo de Temuco"
$ws.Range("C254").Value = "La Araucanía"
$ws.Range("D254").Value = 44726
$ws.Range("E254").Value = 9
$ws.Range("F254").Value = 100112037
$ws.Range("G254").Value = "Cebollín"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 35
$ws.Range("K254").Value = 7000
$ws.Range("L254").Value = 7000
$ws.Range("M254").Value = 7000
$ws.Range("N254").Value = "`$/docena de paquetes"
$ws.Range("O254").Value = "Región del Maule"
$ws.Range("P254").Value = 583
$ws.Range("Q254").Value = 12
$ws.Range("R254").Value = "Hortaliza"

# Apply the same number format as the other date cells in column D.
$ws.Range("D254").NumberFormat = $ws.Range("D255").NumberFormat
